$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 297, shifting existing rows
# 297-312 down to 299-314 (matches dimension growing from R312 to R314).
$ws.Rows("297:298").Insert()

# New row 297: Primera, Volumen 300, Precio [1200,1300,1250]
$ws.Range("A297").Value = 7
$ws.Range("B297").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C297").Value = "Ñuble"
$ws.Range("D297").Value = 44939
$ws.Range("E297").Value = 16
$ws.Range("F297").Value = 100112006
$ws.Range("G297").Value = "Repollo"
$ws.Range("H297").Value = "Crespo record"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 300
$ws.Range("K297").Value = 1200
$ws.Range("L297").Value = 1300
$ws.Range("M297").Value = 1250
$ws.Range("N297").Value = "`$/unidad"
$ws.Range("O297").Value = "Provincia de Diguillín"
$ws.Range("P297").Value = 1250
$ws.Range("Q297").Value = 1
$ws.Range("R297").Value = "Hortaliza"

# New row 298: Segunda, Volumen 300, Precio [1000,1000,1000]
$ws.Range("A298").Value = 7
$ws.Range("B298").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C298").Value = "Ñuble"
$ws.Range("D298").Value = 44939
$ws.Range("E298").Value = 16
$ws.Range("F298").Value = 100112006
$ws.Range("G298").Value = "Repollo"
$ws.Range("H298").Value = "Crespo record"
$ws.Range("I298").Value = "Segunda"
$ws.Range("J298").Value = 300
$ws.Range("K298").Value = 1000
$ws.Range("L298").Value = 1000
$ws.Range("M298").Value = 1000
$ws.Range("N298").Value = "`$/unidad"
$ws.Range("O298").Value = "Provincia de Diguillín"
$ws.Range("P298").Value = 1000
$ws.Range("Q298").Value = 1
$ws.Range("R298").Value = "Hortaliza"
